$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Daily GitHub Actions refresh of cryptocurrency price / Volume(1h) data.
# Rows 29-31 also have a reordering of coins (Fetch.AI/PEPE/Aptos -> Aptos/Fetch.AI/PEPE).
# Price cells (column D) are forced to Text format ("@") before assignment so
# values such as "1.00", "0.0000141" or "60.037.98" are preserved verbatim
# instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.037.98'
$ws.Range("E2").Value = '  -3.76%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.529.95'
$ws.Range("E3").Value = '  -3.56%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '546.15'
$ws.Range("E5").Value = '  -1.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.70'
$ws.Range("E6").Value = '  -4.23%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.82%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.529.17'
$ws.Range("E9").Value = '  -3.58%  '
$ws.Range("E10").Value = '  -3.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.159'
$ws.Range("E11").Value = '  -0.86%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.44'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("E13").Value = '  -1.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.970.43'
$ws.Range("E14").Value = '  -3.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '24.79'
$ws.Range("E15").Value = '  -3.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.047.07'
$ws.Range("E16").Value = '  -3.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000141'
$ws.Range("E17").Value = '  -2.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.510.54'
$ws.Range("E18").Value = '  -4.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.59'
$ws.Range("E19").Value = '  -0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.39'
$ws.Range("E20").Value = '  -3.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.66'
$ws.Range("E21").Value = '  -3.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  -0.11%  '
$ws.Range("E23").Value = '  -4.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.65'
$ws.Range("E24").Value = '  -1.91%  '
$ws.Range("E25").Value = '  -9.51%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.51%  '
$ws.Range("E27").Value = '  -2.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.92'
$ws.Range("E28").Value = '  -1.30%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.99'
$ws.Range("E29").Value = '  -1.48%  '
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.30'
$ws.Range("E30").Value = '  -2.69%  '
$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0₃0797'
$ws.Range("E31").Value = '  -3.58%  '
$ws.Range("E32").Value = '  -3.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '159.31'
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("E35").Value = '  +1.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.93'
$ws.Range("E36").Value = '  -1.45%  '
$ws.Range("E37").Value = '  -3.53%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.74'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.09'
$ws.Range("E39").Value = '  -0.78%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '317.41'
$ws.Range("E40").Value = '  -5.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.76'
$ws.Range("E41").Value = '  -2.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.79'
$ws.Range("E42").Value = '  -2.77%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.836'
$ws.Range("E43").Value = '  -6.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.604'
$ws.Range("E45").Value = '  -0.70%  '
$ws.Range("E46").Value = '  -1.66%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '127.16'
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0533'
$ws.Range("E48").Value = '  -2.58%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0944'
$ws.Range("E49").Value = '  -1.57%  '
$ws.Range("E50").Value = '  -1.96%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.72'
$ws.Range("E51").Value = '  -4.91%  '
